$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.568.49"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "'1.881.58"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'1.029"
$ws.Range("E4").Value = "  +2.54%  "
$ws.Range("D5").Value = "'318.61"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'1.026"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "'0.5161"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "'0.3953"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "'0.08337"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'42.26"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'6.270"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.865.56"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").Value = "'20.44"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "'1.029"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").Value = "'7.266"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "'0.00001112"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'91.60"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'0.06791"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "'1.025"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").Value = "'17.72"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "'5.985"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'28.592.94"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'162.26"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.070.70"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'20.85"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'2.376"
$ws.Range("E29").Value = "  -5.50%  "
$ws.Range("D30").Value = "'127.62"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'0.1055"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").Value = "'1.037"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'5.860"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "'3.671"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("D36").Value = "'0.06514"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'9.155"
$ws.Range("E37").Value = "  -5.80%  "
$ws.Range("D38").Value = "'0.2187"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "'1.254"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.190"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6451"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "'4.989"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'11.22"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'0.6039"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.732"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.95"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").Value = "'1.243"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").Value = "'1.996"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'1.217"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'122.41"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'0.06877"
$ws.Range("E51").Value = "  -0.22%  "

$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
